# "Generate Report for Handoff" - mark b.md as ready for handoff across all
# sheets of the localization status report and record the new handoff
# package / timestamps / warning for the out-of-date handback.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$hoDate = "2016-08-25 04:36:40"

$warning = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5ee609761a429679506052445181a54b2f0b96f4/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4f93c0ad25fbb2d0a37e7099544d6050e219a1c/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: update b.md's zh-cn / de-de status and the latest
# handoff xliff generation date.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status
$overview.Range("G3").Value = $hoDate
$overview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# zh-cn sheet: b.md (row 3) just got a fresh handoff package.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-25 04:36:35"
$zhcn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("P3").Value = $warning
$zhcn.Columns.Item(16).ColumnWidth = 39.2

# ---------------------------------------------------------------------
# de-de sheet: b.md (row 3) just got a fresh handoff package.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = $hoDate
$dede.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Range("P3").Value = $warning
$dede.Columns.Item(16).ColumnWidth = 39.2
